# Repull data, push all data, mean calculation
# Update column F (dSF) values for specific rows to reflect repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 5
    4  = -4
    6  = 4
    9  = -3
    11 = 1
    12 = 1
    14 = -4
    15 = 2
    16 = 1
    21 = -4
    23 = 7
    25 = 1
    26 = -6
    29 = -1
    34 = 4
    36 = 5
    39 = 6
    40 = 5
    41 = 3
    42 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
